# Update Name of Algo
# Apply updated values to result_data_RandomForest.xlsx (Sheet1)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = -6.358699999999991
$ws.Range("B3").Value = 6.10929999999999
$ws.Range("D5").Value = -8.207399999999993
$ws.Range("B14").Value = 9.0921
$ws.Range("B16").Value = 8.975
$ws.Range("D16").Value = -7.870799999999998
$ws.Range("B21").Value = 5.774199999999994
$ws.Range("B23").Value = 5.363300000000001
$ws.Range("B25").Value = 5.955899999999994
